$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.236.41"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3
$ws.Range("D3").Value = "2.267.19"
$ws.Range("E3").Value = "  -1.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'306.49"
$ws.Range("E5").Value = "  -0.67%  "

# Row 6
$ws.Range("D6").Value = "'96.87"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("E7").Value = "  -1.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  -1.57%  "

# Row 10
$ws.Range("D10").Value = "'35.07"
$ws.Range("E10").Value = "  -1.77%  "

# Row 11
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  -2.55%  "

# Row 12
$ws.Range("E12").Value = "  -0.05%  "

# Row 13
$ws.Range("D13").Value = "'6.95"
$ws.Range("E13").Value = "  +2.32%  "

# Row 14
$ws.Range("D14").Value = "2.620.47"
$ws.Range("E14").Value = "  -1.50%  "

# Row 15
$ws.Range("D15").Value = "'14.73"
$ws.Range("E15").Value = "  +0.65%  "

# Row 16
$ws.Range("D16").Value = "2.255.15"
$ws.Range("E16").Value = "  -1.79%  "

# Row 17
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = "  -1.13%  "

# Row 18
$ws.Range("D18").Value = "42.119.73"
$ws.Range("E18").Value = "  -1.19%  "

# Row 19
$ws.Range("D19").Value = "'12.31"
$ws.Range("E19").Value = "  -4.24%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -1.95%  "

# Row 21
$ws.Range("D21").Value = "'6.02"
$ws.Range("E21").Value = "  -0.84%  "

# Row 22
$ws.Range("D22").Value = "'67.85"
$ws.Range("E22").Value = "  -0.86%  "

# Row 23
$ws.Range("D23").Value = "'237.79"
$ws.Range("E23").Value = "  -2.97%  "

# Row 24: 'ImmutableX' -> 'PancakeSwap'
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.58"
$ws.Range("E24").Value = "  -1.88%  "

# Row 25: 'PancakeSwap' -> 'ImmutableX'
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'1.97"
$ws.Range("E25").Value = "  -0.30%  "

# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("D27").Value = "'23.56"
$ws.Range("E27").Value = "  -3.36%  "

# Row 28
$ws.Range("D28").Value = "'37.86"
$ws.Range("E28").Value = "  +1.63%  "

# Row 29
$ws.Range("D29").Value = "'9.59"
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("E30").Value = "  +0.30%  "

# Row 31
$ws.Range("D31").Value = "'163.24"
$ws.Range("E31").Value = "  +0.98%  "

# Row 32
$ws.Range("D32").Value = "'5.25"
$ws.Range("E32").Value = "  -2.77%  "

# Row 33
$ws.Range("E33").Value = "  +0.16%  "

# Row 34
$ws.Range("E34").Value = "  +1.54%  "

# Row 35
$ws.Range("D35").Value = "'17.64"
$ws.Range("E35").Value = "  +1.12%  "

# Row 36
$ws.Range("D36").Value = "'0.0738"
$ws.Range("E36").Value = "  -2.76%  "

# Row 37
$ws.Range("E37").Value = "  -0.81%  "

# Row 38
$ws.Range("E38").Value = "  -4.65%  "

# Row 39
$ws.Range("D39").Value = "'1.82"
$ws.Range("E39").Value = "  -1.52%  "

# Row 40
$ws.Range("E40").Value = "  -1.91%  "

# Row 41
$ws.Range("D41").Value = "'4.06"
$ws.Range("E41").Value = "  -3.93%  "

# Row 42
$ws.Range("E42").Value = "  +2.72%  "

# Row 43
$ws.Range("D43").Value = "'19.12"
$ws.Range("E43").Value = "  -5.02%  "

# Row 44
$ws.Range("D44").Value = "1.953.05"
$ws.Range("E44").Value = "  -3.63%  "

# Row 45
$ws.Range("D45").Value = "'0.0281"
$ws.Range("E45").Value = "  -1.73%  "

# Row 46
$ws.Range("D46").Value = "'2.92"
$ws.Range("E46").Value = "  -3.58%  "

# Row 47
$ws.Range("D47").Value = "'9.86"
$ws.Range("E47").Value = "  -4.30%  "

# Row 48
$ws.Range("D48").Value = "'53.86"
$ws.Range("E48").Value = "  -0.36%  "

# Row 49
$ws.Range("D49").Value = "'92.40"
$ws.Range("E49").Value = "  -1.11%  "

# Row 50
$ws.Range("E50").Value = "  -2.83%  "

# Row 51
$ws.Range("E51").Value = "  -2.69%  "
